# Weekly driver report update for 2025-04-20
#
# 1) Update the "Bad Drivers" summary table (rows 3-6): Critical Minutes and
#    Good Roaming Calculation (%) figures refreshed, plus the Totals row.
# 2) Refresh the "Good Drivers" table (rows 14-27): the rows are re-sorted by
#    Driver Vintage (newest first) with a couple of updated sample counts and
#    several newly-populated Driver Vintage dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Bad Drivers summary table
# ---------------------------------------------------------------------
$ws.Range("C3").Value = 434
$ws.Range("D3").Value = 93.7

$ws.Range("C4").Value = 129
$ws.Range("D4").Value = 97

$ws.Range("C5").Value = 2206
$ws.Range("D5").Value = 98.5

$ws.Range("C6").Value = 2769

# ---------------------------------------------------------------------
# 2) Good Drivers table (rows 14-27) - re-sorted + refreshed values
# ---------------------------------------------------------------------
# Each entry: Adapter-Driver name, Total Samples, Good Roaming (%), Driver Vintage (text date)
$goodDrivers = @(
    @("Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4",      445055, 99.90000000000001, "2024-11-10"),
    @("Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9",        77849, 99.90000000000001, "2021-08-18"),
    @("Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1",        34244, 100,               "2021-04-27"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11", 170510, 99.90000000000001, "2021-01-19"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1",    52096, 100,               "2020-09-28"),
    @("Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2",        59673, 100,               "2020-08-05"),
    @("Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6",        113652, 100,               "2020-01-06"),
    @("Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1",         56018, 100,               "2019-12-14"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3", 161874, 100,               "2019-09-05"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5", 143342, 99.90000000000001, "2019-08-25"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2",  20227, 100,               "2019-05-11"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1",   34065, 100,               "2019-04-28"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1",   48540, 100,               "2019-03-16"),
    @("Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2",  184564, 99.90000000000001, "2018-11-25")
)

$row = 14
foreach ($driver in $goodDrivers) {
    $ws.Range("A$row").Value = $driver[0]
    $ws.Range("B$row").Value = $driver[1]
    $ws.Range("D$row").Value = $driver[2]

    # Force column E to be stored as plain text so the vintage date ("2024-11-10")
    # is not silently reinterpreted as a date serial number.
    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $driver[3]

    $row = $row + 1
}
